$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header row: column F header text changes from "Email" to "Email Address".
#    (Column G keeps its "Email Subject" text.)
$ws.Range("F1").Value = "Email Address"

# 2. Mark F2 with a (cosmetically no-op) fill flag, matching the authored edit
#    that re-applied "No Fill" to F2 via the Format Cells / Fill Color UI.
$ws.Range("F2").Interior.ColorIndex = -4142   # xlColorIndexNone

# 3. Append the new rows of partnership-email data (rows 9-14), following the
#    same repeating 3-row pattern (Tech-Neo / GlobalTech / InnovateNow) already
#    present for the earlier batches, but with the two new "Date Processed" timestamps.

$pattern = @(
    @{ B = "Moris Mwai"; C = "Tech-Neo GmbH";             D = "Am main City, Germany`n";             E = "DE1567890";   F = "morismwai1@gmail.com"; G = "Partnership Offer" },
    @{ B = "Moris Mwai"; C = "GlobalTech Solutions GmbH"; D = "Hauptstrasse 24, Berlin, Germany`n";   E = "DE789654321"; F = "morismwai1@gmail.com"; G = "Partnership Offer" },
    @{ B = "Moris Mwai"; C = "InnovateNow GmbH";           D = "Bahnhofstrasse 45, Munich, Germany`n"; E = "DE345678901"; F = "morismwai1@gmail.com"; G = "Join our company as an investor" }
)

$dates = @(45859.800729166665, 45859.80741898148)

$r = 9
foreach ($d in $dates) {
    foreach ($row in $pattern) {
        $ws.Range("A$r").Value = $d
        $ws.Range("A$r").NumberFormat = "m/d/yy h:mm"
        $ws.Range("B$r").Value = $row.B
        $ws.Range("C$r").Value = $row.C
        $ws.Range("D$r").Value = $row.D
        $ws.Range("E$r").Value = $row.E
        $ws.Range("F$r").Value = $row.F
        $ws.Range("G$r").Value = $row.G
        # Writing a value that embeds a line break auto-expands the row's
        # height; AutoFit restores the row back to the sheet's normal
        # (non-custom) height, matching the source rows' plain formatting.
        $ws.Rows.Item($r).AutoFit()
        $r = $r + 1
    }
}

# 4. Update the active selection to F2, matching the authored edit.
$ws.Range("F2").Select()
